# Add Ukrainian translations to the "Translations" column of both sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Multilingual_spreadsheet_sample")
$ws2 = $wb.Worksheets.Item("Sheet1")

# Sheet "Multilingual_spreadsheet_sample" - Translations live in column E
$ws1.Range("E2").Value = "Переклад 1"
$ws1.Range("E3").Value = "Переклад 2"
$ws1.Range("E4").Value = "Переклад 3"

# Sheet "Sheet1" - Translations live in column E as well
$ws2.Range("E2").Value = "Переклад 1"
$ws2.Range("E3").Value = "Переклад 2"
$ws2.Range("E4").Value = "Переклад 3"
